# Scenario-mode "Remanente" column addition to the monthly program sheet.
# Adds a new column Q ("Remanente 31-01-2026") and re-points the P column
# formula from +O*N to +O-Q, filling in the new Q values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Programa Enero")
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Extend formatting from column P into the new column Q so the new
#    cells inherit the same look (borders/number format/font) as the
#    rest of the table, row by row.
# ---------------------------------------------------------------------

# Row 1: plain placeholder cell with no style (mirrors the other
# unstyled placeholder cells used elsewhere on the sheet).
$ws.Range("P50").Copy()
$ws.Range("Q1").PasteSpecial(-4122)

# Header row 2 — copy P2's header style, then set the new header text.
$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)
$ws.Range("Q2").Value2 = "Remanente 31-01-2026"

# "Section title" rows use style 7 (same as the existing O/P columns on
# those rows).
$headerStyleRows = 3,4,5,6,7,10,21,25,27,30,34
$ws.Range("P3").Copy()
foreach ($r in $headerStyleRows) {
    $ws.Range("Q$r").PasteSpecial(-4122)
}

# Data rows use style 2 (same as the existing O/P columns on those
# rows).
$dataStyleRows = 8,9,11,12,13,14,15,16,17,18,19,20,22,23,24,26,28,29,31,32,33,35,36,37,38,39,40,41,42,43,44,45,46,47,48
$ws.Range("P8").Copy()
foreach ($r in $dataStyleRows) {
    $ws.Range("Q$r").PasteSpecial(-4122)
}

# Trailing filter helper rows below the table (empty placeholder cells
# that just extend the row's used range into column Q).
$trailingRows = 50,52,53,54,55,56
$ws.Range("P50").Copy()
foreach ($r in $trailingRows) {
    $ws.Range("Q$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Fill in the "Remanente" values the user measured for January.
# ---------------------------------------------------------------------

$remanente = @{
    8  = 156.03999999999996
    9  = 19.5
    11 = 23.75
    12 = 64.245999999999995
    13 = 69.779999999999973
    14 = 15.738
    15 = 4.9984999999999999
    16 = 204.99600000000001
    17 = 60
    18 = 116.4
    19 = 50
    20 = 624.976
    22 = 78.499000000000024
    23 = 104.05500000000001
    24 = 23.033000000000015
    26 = 189.59999999999991
    28 = 21.009000000000015
    29 = 51.003999999999991
    31 = 110.40000000000009
    32 = 190
    33 = 890
    35 = 255.303
    36 = 26
    37 = 190.48999999999998
    38 = 3.6539999999999964
    39 = 135
    40 = 144.50640000000001
    41 = 64.752000000000066
    42 = 70
    43 = 221.00650000000002
    44 = 8.724899999999991
    45 = 5.8500000000000085
    46 = 8.3693999999999988
    47 = 2.1899999999999977
    48 = 18
}

foreach ($r in $remanente.Keys) {
    $ws.Range("Q$r").Value2 = $remanente[$r]
}

# ---------------------------------------------------------------------
# 3) Re-point column P so it is the "liberado" quantity minus the
#    remanente instead of O*N.
# ---------------------------------------------------------------------

$formulaRows = 8,9,11,12,13,14,15,16,17,18,19,20,22,23,24,26,28,29,31,32,33,35,36,37,38,39,40,41,42,43,44,45,46,47,48
foreach ($r in $formulaRows) {
    $ws.Range("P$r").Formula = "=+O$r-Q$r"
}

# ---------------------------------------------------------------------
# 4) Restore the selection to match where the user was working.
# ---------------------------------------------------------------------

$ws.Range("Q35:Q48").Select()

$wb.Application.Calculate()
